$wb = $excel.ActiveWorkbook

# Add the new "CMS" sheet as the last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cms = $wb.Worksheets.Add($null, $lastSheet)
$cms.Name = "CMS"

# Header row values (Contact_* / OM_* extract columns).
$headers = @(
    "Contact_ID`n",
    "Contact_Date`n",
    "Contact_Type_Code",
    "Contact_Type_Desc",
    "Contact_Staff_Name",
    "Contact_Staff_Key",
    "Contact_Staff_Grade",
    "Contact_Team_Key",
    "Contact_Provider_Code",
    "OM_Name`n",
    "OM_Key`n",
    "OM_Grade`n",
    "OM_Team_Key`n",
    "OM_Provider_Code`n"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cms.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Format the header row: 9pt white-on-... Arial (#333333), solid white fill,
# left-aligned, with a taller custom row height.
$headerRange = $cms.Range("A1:N1")
$headerRange.Font.Size = 9
$headerRange.Font.Name = "Arial"
$headerRange.Font.Color = 3355443
$headerRange.Interior.Pattern = 1
$headerRange.Interior.Color = 16777215
$headerRange.Interior.PatternColor = 16777215
$headerRange.HorizontalAlignment = -4131
$cms.Rows.Item(1).RowHeight = 23.25

# Select the whole header row (A1:XFD1) and make CMS the active/selected tab.
$cms.Range("A1:XFD1").Select() | Out-Null
$cms.Activate() | Out-Null
